$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices and 1h volume deltas),
# plus a couple of row-content swaps, per the Mar 23 2024 refresh.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.417.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.348.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.20%  "

# Row 7
$ws.Range("E7").Value = "  -2.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.340.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.40%  "

# Row 9
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.627"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.53%  "

# Row 11
$ws.Range("E11").Value = "  +0.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.42%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.884.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.41%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.29%  "

# Row 17
$ws.Range("E17").Value = "  -2.88%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.348.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.326.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.979"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "435.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.74%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.90%  "

# Row 24
$ws.Range("E24").Value = "  -4.87%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "84.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27
$ws.Range("E27").Value = "  -2.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.63%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "577.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.28%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.79%  "

# Row 36
$ws.Range("E36").Value = "  +0.11%  "

# Row 37
$ws.Range("E37").Value = "  -8.13%  "

# Row 38
$ws.Range("E38").Value = "  -3.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.91%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.367"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.58%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.101.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.20%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.81%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0409"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.34%  "

# Row 48
$ws.Range("E48").Value = "  -1.91%  "

# Row 49
$ws.Range("E49").Value = "  -3.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.44%  "
